# "integrated gyro + motion" — remove the now-obsolete I2C/gpio check and
# "Run ./dist/start.js" sections at the end of the document, while keeping
# the trailing "_GoBack" bookmark alive on the paragraph that remains.

$d = $word.ActiveDocument

# The bookmark (_GoBack) currently lives inside the "gpio readall" paragraph
# that is about to be deleted. Re-anchor it on the empty paragraph that sits
# right before the "Check I2C bus data:" heading so it survives the cut.
$anchorPara = $d.Paragraphs.Item(9)
$d.Bookmarks.Add("_GoBack", $anchorPara.Range)

# Remove everything from the "Check I2C bus data:" heading through to the
# final trailing empty paragraph (the whole I2C-check + "Run" block).
$firstDoomed = $d.Paragraphs.Item(10)
$lastDoomed = $d.Paragraphs.Item($d.Paragraphs.Count)
$killRange = $d.Range($firstDoomed.Range.Start, $lastDoomed.Range.End)
$killRange.Delete()
